# Apply updated "dSF" (column F) values for hendricks_kyle.xlsx
# Commit message: "repull data, push all data, mean calculation"
# The diff shows column F (dSF) values recalculated for most rows (2-36),
# leaving a few rows (5, 9, 19, 21) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = -1
    4  = 1
    6  = -2
    7  = -3
    8  = 1
    10 = -2
    11 = 3
    12 = -6
    13 = -3
    14 = 3
    15 = -2
    16 = -2
    17 = -1
    18 = 3
    20 = -4
    22 = -6
    23 = -3
    24 = -6
    25 = 1
    26 = -3
    27 = 1
    28 = -3
    29 = 2
    30 = -3
    31 = -1
    32 = -3
    33 = 5
    34 = -2
    35 = -3
    36 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
